$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("frais divers")

$ws.Range("A19").Value = "Test item added programmatically"
$ws.Range("B19").Value = 123.45
$ws.Range("C19").Value = "Entretien"
